# Weekly update: a new price record is inserted at the top of the data
# block (row 258) and all subsequent records shift down by one row, with
# the former last row (310) duplicating into a brand-new row 311.
#
# Columns A,B,C,E,F,G,H,N,Q,R are constant for every data row in this
# block and do not need to be touched. Only D,I,J,K,L,M,O,P differ
# between records, so we shift just those down one row at a time,
# working from the bottom (row 311) up to (row 260) so that each read
# happens before its source row gets overwritten.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colIndex = @{ "D" = 4; "I" = 9; "J" = 10; "K" = 11; "L" = 12; "M" = 13; "O" = 15; "P" = 16 }
$cols = @("D", "I", "J", "K", "L", "M", "O", "P")

for ($r = 311; $r -ge 260; $r--) {
    $src = $r - 1
    foreach ($col in $cols) {
        $c = $colIndex[$col]
        $val = $ws.Cells.Item($src, $c).Value2
        $ws.Cells.Item($r, $c).Value2 = $val
    }
}

# Row 258 becomes the new weekly record. Row 259 already holds the data
# that used to be in row 258 (they were identical before the edit), so it
# is left untouched.
$ws.Cells.Item(258, 4).Value2 = 44504    # D258 Fecha
$ws.Cells.Item(258, 9).Value2 = "Primera" # I258 Calidad (unchanged)
$ws.Cells.Item(258, 10).Value2 = 840      # J258 Volumen
$ws.Cells.Item(258, 11).Value2 = 750      # K258 Precio minimo
$ws.Cells.Item(258, 12).Value2 = 800      # L258 Precio maximo
$ws.Cells.Item(258, 13).Value2 = 773      # M258 Precio promedio ponderado
$ws.Cells.Item(258, 15).Value2 = "Región Metropolitana" # O258 Origen (unchanged)
$ws.Cells.Item(258, 16).Value2 = 258      # P258 Precio $/Kg

# Row 311 is a brand-new row, so the columns that are constant across the
# whole data block (they never differ from one record to the next) must
# be populated explicitly too.
$ws.Cells.Item(311, 1).Value2 = 6                                           # A311 Mercado ID
$ws.Cells.Item(311, 2).Value2 = "Mercado Mayorista Lo Valledor de Santiago" # B311 Mercado
$ws.Cells.Item(311, 3).Value2 = "Metropolitana"                             # C311 Región
$ws.Cells.Item(311, 5).Value2 = 13                                          # E311 Codreg
$ws.Cells.Item(311, 6).Value2 = 100112039                                   # F311 Categoría ID
$ws.Cells.Item(311, 7).Value2 = "Ciboulette"                                # G311 Categoría
$ws.Cells.Item(311, 8).Value2 = "Sin especificar"                           # H311 Variedad
$ws.Cells.Item(311, 14).Value2 = "$/docena de atados"                       # N311 Unidad de comercialización
$ws.Cells.Item(311, 17).Value2 = 3                                          # Q311 Kg o Unidades
$ws.Cells.Item(311, 18).Value2 = "Hortaliza"                                # R311 Clasificación

# New row 311 needs the same number format that column D carries for
# every data row (numFmtId 165 date/time format), matching the rest of
# the block.
$ws.Range("D311").NumberFormat = $ws.Range("D310").NumberFormat
